# Update the adult TOD raw->ss lookup tables on each age-band sheet.
# A new "raw score 0" row is inserted at the top of each table (rows shift
# down by one) and the standard-score (column B) values are refreshed.
$wb = $excel.ActiveWorkbook

# Sheet 1: "18.0-23.11"
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 50
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 51
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 52
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 53
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 54
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 55
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 56
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 57
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 59
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 60
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 61
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 62
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 64
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 65
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 67
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 69
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 70
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 72
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 75
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 77
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 79
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 82
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 86
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 90
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 96
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 111

# Sheet 2: "24.0-39.11"
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 50
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 51
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 52
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 53
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 55
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 56
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 57
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 58
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 59
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 60
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 62
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 63
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 64
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 66
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 68
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 69
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 71
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 73
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 75
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 77
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 80
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 83
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 86
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 91
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 97
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 111

# Sheet 3: "40.0-49.11"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 52
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 53
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 55
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 56
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 57
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 58
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 59
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 60
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 62
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 63
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 64
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 66
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 67
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 69
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 71
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 72
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 74
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 76
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 79
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 81
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 84
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 88
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 92
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 98
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 111

# Sheet 4: "50.0-59.11"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 53
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 54
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 55
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 56
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 57
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 58
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 59
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 60
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 62
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 63
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 64
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 66
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 67
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 69
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 70
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 72
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 74
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 76
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 78
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 80
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 83
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 86
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 89
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 94
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 99
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 112

# Sheet 5: "60.0-69.11"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 55
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 56
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 57
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 58
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 59
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 60
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 61
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 63
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 64
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 65
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 66
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 68
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 69
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 71
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 73
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 74
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 76
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 78
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 80
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 83
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 85
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 88
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 91
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 95
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 101
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 112

# Sheet 6: "70.0-89.11"
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 59
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 60
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 61
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 62
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 63
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 64
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 66
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 67
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 68
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 69
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 71
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 72
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 74
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 75
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 77
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 79
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 80
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 82
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 84
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 87
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 89
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 92
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 95
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 98
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 103
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 113
